$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old ShiftCode value in row 2 (its only shared-string use), so the
# unused "AA" string is dropped before new strings are introduced.
$ws.Range("E2").ClearContents()

# --- EmpUnqID (A) for the employee's sanction rows ---
$ws.Range("A2").Value = 100234
$ws.Range("A3").Value = 100234
$ws.Range("A4").Value = 100234
$ws.Range("A5").Value = 100234

# --- SanDate (B) : one row per day the overtime sanction applies to ---
$ws.Range("B2").Value = "'2018-04-04"

# --- InTime (C) for the first date ---
$ws.Range("C2").NumberFormat = "h:mm"
$ws.Range("C2").Value = "'08:00"

$ws.Range("B3").Value = "'2018-04-05"
$ws.Range("B4").Value = "'2018-04-06"
$ws.Range("B5").Value = "'2018-04-07"

# --- OutTime (D) for the first date ---
$ws.Range("D2").NumberFormat = "h:mm"
$ws.Range("D2").Value = "'20:00"

# --- Fill remaining InTime / OutTime cells with the same shift times ---
$ws.Range("C3").NumberFormat = "h:mm"
$ws.Range("C3").Value = "'08:00"
$ws.Range("D3").NumberFormat = "h:mm"
$ws.Range("D3").Value = "'20:00"

$ws.Range("C4").NumberFormat = "h:mm"
$ws.Range("C4").Value = "'08:00"
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("D4").Value = "'20:00"

$ws.Range("C5").NumberFormat = "h:mm"
$ws.Range("C5").Value = "'08:00"
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("D5").Value = "'20:00"

# TPAHours for the last sanction row
$ws.Range("F5").Value = 0

# --- Selection matches the saved cursor position in the edited file ---
$ws.Range("D6").Select()
